# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert 10 blank rows after row 20 (current last data row), copying
#    the "normal" row (row 19) formatting into each inserted row so the
#    new rows 21..30 share formatting with the existing body rows.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 10; $i++) {
    $ws.Range("B19:J19").Copy()
    $ws.Range("B21:J21").EntireRow.Insert()
}

# ---------------------------------------------------------------------
# 2) Row 20 currently still carries the special "last row" formatting
#    (bottom border). Copy that formatting down onto the new last row
#    (row 30) before we overwrite row 20 with normal formatting.
# ---------------------------------------------------------------------
$ws.Range("B20:J20").Copy()
$ws.Range("B30:J30").PasteSpecial(-4122)

# Now make row 20 use the regular body-row formatting (copied from row 19).
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Write the new worker/period data into rows 16..30.
# ---------------------------------------------------------------------
$data = @(
    @("CC", "73193703",   "EDGARD DE JESUS BARRIOS BUENDIA",     "2102", 96000, 2400000),
    @("CC", "1128049639", "DARWIN JAVIER MARTINEZ CASTILLA",     "1610", 56000, 1400000),
    @("CC", "1128049639", "DARWIN JAVIER MARTINEZ CASTILLA",     "1609", 56000, 1400000),
    @("CC", "1013601861", "DAVID FELIPE PINTO ACOSTA",           "2205", 33333, 1000000),
    @("CC", "1143408511", "CARLOS JAVIER LOZANO FRANCO",         "2205", 43181, 1245600),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2507", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2506", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2505", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2504", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2503", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2502", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2501", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2412", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2411", 52000, 1689768),
    @("CC", "1047466573", "ANDREA CAROLINA MARTINEZ BALLESTAS",  "2410", 52000, 1689768)
)

$row = 16
foreach ($rec in $data) {
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 4) Update the summary header cells.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 804514
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 14

# ---------------------------------------------------------------------
# 5) Widen column D to fit the longest worker name.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).EntireColumn.AutoFit()

"done"
